$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RB")

# Add new player row (Week 13 logging) for A.Peterson
$ws.Range("A7").Value = "A.Peterson"
$ws.Range("B7:J7").Value = 0

# Make RB the active sheet/tab and set its selection
$ws.Activate()
$ws.Range("J8").Select()
